$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = 9910.1200000000008
$ws.Range("B4").Value = 9927
$ws.Range("C4").Value = 309.02999999999997
$ws.Range("D4").Value = 309.55
$ws.Range("E4").Value = $true
$ws.Range("F4").Value = 0.17
$ws.Range("G4").Value = 42608.624768518515
$ws.Range("H4").Value = $false

# Row 5
$ws.Range("A5").Value = 9893.27
$ws.Range("B5").Value = 9910.1200000000008
$ws.Range("C5").Value = 309.02999999999997
$ws.Range("D5").Value = 309.55
$ws.Range("E5").Value = $true
$ws.Range("F5").Value = 0.17
$ws.Range("G5").Value = 42608.63784722222
$ws.Range("H5").Value = $false

# Copy the number format from G3 (date style) to G4:G5 to match style index 1
$ws.Range("G3").Copy()
$ws.Range("G4:G5").PasteSpecial(-4122) # xlPasteFormats
